# Strip the stray trailing "16" that got appended to every verse reference
# in column A (e.g. "Matthew 1:316" -> "Matthew 1:3"), for human readability.
# Row 1 is the "Reference"/"Text" header and is left untouched; data runs
# from row 2 through row 452.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 452; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $v = $cell.Value()
    if ($v -ne $null -and $v.EndsWith("16")) {
        $cell.Value = $v.Substring(0, $v.Length - 2)
    }
}
